$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Hardest Puzzle" sheet: selection moves from K12 to D2
# ---------------------------------------------------------------------------
$wsHardest = $wb.Worksheets.Item("Hardest Puzzle")
$wsHardest.Activate()
$wsHardest.Range("D2").Select()

# ---------------------------------------------------------------------------
# 2. "Puzzles from websudoku - level1" sheet: selection moves from H13 to H11
#    (it also loses tabSelected, which happens automatically once another
#    sheet becomes active later in the script)
# ---------------------------------------------------------------------------
$wsLevel1 = $wb.Worksheets.Item("Puzzles from websudoku - level1")
$wsLevel1.Activate()
$wsLevel1.Range("H11").Select()

# ---------------------------------------------------------------------------
# 3. New sheet "Solve2" appended at the end
# ---------------------------------------------------------------------------
$solve2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$solve2.Name = "Solve2"

$solve2.Range("A1").Value = 57719598898
$solve2.Range("B1").Formula = "=A1/1000000000"
$solve2.Range("A2").Value = 58023811826
$solve2.Range("B2").Formula = "=A2/1000000000"

$solve2.Columns("A").ColumnWidth = 17.666666666666668
$solve2.Columns("B").ColumnWidth = 19.833333333333332

$solve2.Activate()
$solve2.Range("C2").Select()

# ---------------------------------------------------------------------------
# 4. New sheet "Solve3" appended at the end
# ---------------------------------------------------------------------------
$solve3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$solve3.Name = "Solve3"

$solve3.Range("A1").Value = "Async"
$solve3.Range("D1").Value = "Seq"

$solve3.Range("A2").Value = 37647574863
$solve3.Range("B2").Formula = "=A2/1000000000"
$solve3.Range("D2").Value = 84320732746
$solve3.Range("E2").Formula = "=D2/1000000000"

$solve3.Range("A3").Value = 35351821976
$solve3.Range("B3").Formula = "=A3/1000000000"
$solve3.Range("D3").Value = 82711369779
$solve3.Range("E3:E19").Formula = "=D3/1000000000"
$solve3.Range("E4").ClearContents()

$solve3.Range("A5").Value = 45335242553
$solve3.Range("A6").Value = 43258363231
$solve3.Range("A8").Value = 42267225511
$solve3.Range("A9").Value = 29921752119
$solve3.Range("B5:B21").Formula = "=A5/1000000000"

$solve3.Columns("A").ColumnWidth = 17.666666666666668
$solve3.Columns("B").ColumnWidth = 21.333333333333332
$solve3.Columns("D").ColumnWidth = 15.166666666666666
$solve3.Columns("E").ColumnWidth = 13.166666666666666

$solve3.Activate()
$solve3.Range("B9").Select()

# ---------------------------------------------------------------------------
# 5. New sheet "Solve4" appended at the end (ends up the active tab)
# ---------------------------------------------------------------------------
$solve4 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$solve4.Name = "Solve4"

$solve4.Range("A1").Value = 117717162966
$solve4.Range("B1").Formula = "=A1/1000000000"

$solve4.Range("A2").Value = 82244142373
$solve4.Range("A5").Value = 38497570580
$solve4.Range("A6").Value = 39724303789
$solve4.Range("A7").Value = 38246429884
$solve4.Range("A9").Value = 23403074462
$solve4.Range("A10").Value = 27738861755
$solve4.Range("A12").Value = 34065698346
$solve4.Range("B2:B21").Formula = "=A2/1000000000"

$solve4.Columns("A").ColumnWidth = 19.833333333333332
$solve4.Columns("B").ColumnWidth = 10.333333333333334

$solve4.Activate()
$solve4.Range("D12").Select()
